# Auto-generated script applying Goblin_Profits.xlsx value updates
# (scheduled-runner refresh of market board price data)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 10005.5
$ws.Range("I7").Value = 5
$ws.Range("K7").Value = 5
$ws.Range("M7").Value = 107
$ws.Range("H14").Value = 10005.5
$ws.Range("I14").Value = 5
$ws.Range("K14").Value = 5
$ws.Range("M14").Value = 186
$ws.Range("H15").Value = 1294.2449
$ws.Range("I15").Value = 1294.2449
$ws.Range("K15").Value = 3882.7347
$ws.Range("M15").Value = -3713.7347
$ws.Range("H138").Value = 1485899.8
$ws.Range("I138").Value = 6960.1
$ws.Range("K138").Value = 20880.3
$ws.Range("M138").Value = -15740.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H38").Value = 16679.75
$ws.Range("I38").Value = 14832.667
$ws.Range("J38").Value = 22221
$ws.Range("K38").Value = 14832.667
$ws.Range("L38").Value = 22221
$ws.Range("M38").Value = -14365.667
$ws.Range("N38").Value = -23155
$ws.Range("H74").Value = 2487.138
$ws.Range("I74").Value = 2065.6191
$ws.Range("K74").Value = 2065.6191
$ws.Range("M74").Value = -1191.6191
$ws.Range("H77").Value = 2487.138
$ws.Range("I77").Value = 2065.6191
$ws.Range("K77").Value = 10328.0955
$ws.Range("M77").Value = -5960.095499999999
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H132").Value = 1679.3043
$ws.Range("I132").Value = 1710.1818
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 5130.5454
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -2600.5454
$ws.Range("N132").Value = -8060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H109").Value = 50000
$ws.Range("J109").Value = 50000
$ws.Range("L109").Value = 50000
$ws.Range("N109").Value = -52774

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 183.9
$ws.Range("J7").Value = 99.833336
$ws.Range("L7").Value = 99.833336
$ws.Range("N7").Value = -325.833336
$ws.Range("H22").Value = 1267.6666
$ws.Range("I22").Value = 659.4286
$ws.Range("K22").Value = 659.4286
$ws.Range("M22").Value = -309.4286
$ws.Range("H31").Value = 3813.3
$ws.Range("I31").Value = 2371.5293
$ws.Range("J31").Value = 5698.6924
$ws.Range("K31").Value = 2371.5293
$ws.Range("L31").Value = 5698.6924
$ws.Range("M31").Value = -2076.5293
$ws.Range("N31").Value = -6288.6924
$ws.Range("H34").Value = 3813.3
$ws.Range("I34").Value = 2371.5293
$ws.Range("J34").Value = 5698.6924
$ws.Range("K34").Value = 2371.5293
$ws.Range("L34").Value = 5698.6924
$ws.Range("M34").Value = -2169.5293
$ws.Range("N34").Value = -6102.6924
$ws.Range("H58").Value = 3130.45
$ws.Range("I58").Value = 3083.8333
$ws.Range("K58").Value = 3083.8333
$ws.Range("M58").Value = -2880.8333
$ws.Range("H132").Value = 1114.5416
$ws.Range("I132").Value = 1092.8572
$ws.Range("K132").Value = 3278.5716
$ws.Range("M132").Value = -748.5715999999998
$ws.Range("H136").Value = 3130.45
$ws.Range("I136").Value = 3083.8333
$ws.Range("K136").Value = 9251.499899999999
$ws.Range("M136").Value = -6701.499899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1859.6428
$ws.Range("I5").Value = 1295.125
$ws.Range("K5").Value = 3885.375
$ws.Range("M5").Value = -3773.375
$ws.Range("H21").Value = 485.6154
$ws.Range("I21").Value = 520.8333
$ws.Range("K21").Value = 1562.4999
$ws.Range("M21").Value = -1389.4999
$ws.Range("H29").Value = 233.3077
$ws.Range("J29").Value = 272.5
$ws.Range("L29").Value = 817.5
$ws.Range("N29").Value = -1371.5
$ws.Range("H127").Value = 2477.3333
$ws.Range("J127").Value = 2477.3333
$ws.Range("L127").Value = 7431.999899999999
$ws.Range("N127").Value = -17351.9999
$ws.Range("H135").Value = 1859.6428
$ws.Range("I135").Value = 1295.125
$ws.Range("K135").Value = 11656.125
$ws.Range("M135").Value = -9121.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 1801.3334
$ws.Range("J6").Value = 2198
$ws.Range("L6").Value = 2198
$ws.Range("N6").Value = -2424
$ws.Range("H16").Value = 1801.3334
$ws.Range("J16").Value = 2198
$ws.Range("L16").Value = 2198
$ws.Range("N16").Value = -2698
$ws.Range("H114").Value = 40000
$ws.Range("J114").Value = 40000
$ws.Range("L114").Value = 40000
$ws.Range("N114").Value = -48678
$ws.Range("H122").Value = 5318.037
$ws.Range("I122").Value = 3971.5454
$ws.Range("K122").Value = 11914.6362
$ws.Range("M122").Value = -9464.636200000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 136011.3
$ws.Range("I2").Value = 146988.16
$ws.Range("J2").Value = 630
$ws.Range("K2").Value = 146988.16
$ws.Range("L2").Value = 630
$ws.Range("M2").Value = -146876.16
$ws.Range("N2").Value = -854
$ws.Range("H7").Value = 8363.333000000001
$ws.Range("I7").Value = 8295
$ws.Range("J7").Value = 8500
$ws.Range("K7").Value = 8295
$ws.Range("L7").Value = 8500
$ws.Range("M7").Value = -8183
$ws.Range("N7").Value = -8724
$ws.Range("H22").Value = 3991.25
$ws.Range("J22").Value = 4998.5
$ws.Range("L22").Value = 4998.5
$ws.Range("N22").Value = -5588.5
$ws.Range("H27").Value = 3991.25
$ws.Range("J27").Value = 4998.5
$ws.Range("L27").Value = 4998.5
$ws.Range("N27").Value = -5212.5
$ws.Range("H61").Value = 2868.3794
$ws.Range("I61").Value = 1166.619
$ws.Range("J61").Value = 7335.5
$ws.Range("K61").Value = 1166.619
$ws.Range("L61").Value = 7335.5
$ws.Range("M61").Value = -964.6189999999999
$ws.Range("N61").Value = -7739.5
$ws.Range("H103").Value = 23110.4
$ws.Range("J103").Value = 23110.4
$ws.Range("L103").Value = 23110.4
$ws.Range("N103").Value = -25454.4
$ws.Range("H113").Value = 2868.3794
$ws.Range("I113").Value = 1166.619
$ws.Range("J113").Value = 7335.5
$ws.Range("K113").Value = 1166.619
$ws.Range("L113").Value = 7335.5
$ws.Range("M113").Value = 1003.381
$ws.Range("N113").Value = -11675.5
$ws.Range("H126").Value = 8363.333000000001
$ws.Range("I126").Value = 8295
$ws.Range("J126").Value = 8500
$ws.Range("K126").Value = 24885
$ws.Range("L126").Value = 25500
$ws.Range("M126").Value = -22415
$ws.Range("N126").Value = -30440

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 433.48
$ws.Range("I107").Value = 401.4
$ws.Range("J107").Value = 561.8
$ws.Range("K107").Value = 1204.2
$ws.Range("L107").Value = 1685.4
$ws.Range("M107").Value = 715.8000000000002
$ws.Range("N107").Value = -5525.4

Write-Host "Applied 170 cell updates across 8 sheets"